$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: delivery date "02/06" -> "27/10"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("02/06", $true, $false, $false, $false, $false, $true, 1, $false, "27/10", 2)

# ---------------------------------------------------------------------------
# Change 2: "...grupos de 3 pessoas." -> "...grupos de 3 ou 4 pessoas."
#           and move the "_GoBack" bookmark from after "O que entregar?" to
#           right after the new " ou 4" text (before the new " pessoas" run).
# ---------------------------------------------------------------------------

# Remove the existing _GoBack bookmark (currently right after "O que entregar?")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the run that currently holds " pessoas" (right after the "3" run)
$full = $d.Content.Text
$pos = $full.IndexOf(" pessoas")

# Protect the boundary between "3" and " pessoas" with a temporary bookmark so
# that editing the text in place does not get merged backwards into the "3" run.
$bmRangeL = $d.Range($pos, $pos)
$d.Bookmarks.Add("TEMPLEFT", $bmRangeL)

# Change the text of that run from " pessoas" to " ou 4"
$full1 = $d.Content.Text
$pos1 = $full1.IndexOf(" pessoas")
$r = $d.Range($pos1, $pos1 + 8)
$r.Text = " ou 4"

# Insert a one-character placeholder right after " ou 4" - this will become the
# new " pessoas" run. Using a placeholder lets us sandwich it between two
# bookmarks before giving it its final text/formatting, which prevents the
# engine from merging it into its same-formatted neighbours.
$full2 = $d.Content.Text
$posAfter = $full2.IndexOf(" ou 4") + 5
$ph = $d.Range($posAfter, $posAfter)
$ph.InsertAfter("X")

# Bookmark right after the placeholder (protects the right boundary, against the
# following ". " run)
$full3 = $d.Content.Text
$posAfterX = $full3.IndexOf(" ou 4X") + 6
$bmRangeR = $d.Range($posAfterX, $posAfterX)
$d.Bookmarks.Add("TEMPRIGHT", $bmRangeR)

# Bookmark right before the placeholder: this is exactly where "_GoBack" needs
# to end up (between " ou 4" and the new " pessoas" run), so create it with its
# final name directly.
$full4 = $d.Content.Text
$posBeforeX = $full4.IndexOf("X")
$bmRangeGoBack = $d.Range($posBeforeX, $posBeforeX)
$d.Bookmarks.Add("_GoBack", $bmRangeGoBack)

# Replace the placeholder with the final text and apply bold/italic/blue formatting
$full5 = $d.Content.Text
$posX = $full5.IndexOf("X")
$rX = $d.Range($posX, $posX + 1)
$rX.Text = " pessoas"
$newRunRange = $d.Range($posX, $posX + 8)
$newRunRange.Font.Bold = $true
$newRunRange.Font.Italic = $true
$newRunRange.Font.Color = 12611584 ## RGB(0x00,0x70,0xC0) == wdColor 0070C0
$newRunRange.Font.Size = 16

# Remove the temporary protection bookmarks, keeping only "_GoBack"
$d.Bookmarks("TEMPLEFT").Delete()
$d.Bookmarks("TEMPRIGHT").Delete()

Write-Host "Edits applied successfully"
